$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.241.28"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.71%  '

$ws.Range('D3').Value = "'1.657.89"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.43%  '

$ws.Range('E4').Value = '  +0.53%  '

$ws.Range('D5').Value = "'218.33"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.13%  '

$ws.Range('D6').Value = "'0.5320"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.38%  '

$ws.Range('E7').Value = '  +0.44%  '

$ws.Range('D8').Value = "'0.2628"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.64%  '

$ws.Range('D9').Value = "'0.06351"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.97%  '

$ws.Range('D10').Value = "'20.46"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.77%  '

$ws.Range('D11').Value = "'0.07839"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.21%  '

$ws.Range('D12').Value = "'4.538"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.48%  '

$ws.Range('D13').Value = "'1.652.98"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.07%  '

$ws.Range('D14').Value = "'1.884.70"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.40%  '

$ws.Range('D15').Value = "'0.5513"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.30%  '

$ws.Range('D16').Value = "'0.0₅8185"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.05%  '

$ws.Range('D17').Value = "'65.52"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.59%  '

$ws.Range('D18').Value = "'26.215.32"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.49%  '

$ws.Range('E19').Value = '  +0.57%  '

$ws.Range('D20').Value = "'4.620"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.45%  '

$ws.Range('D21').Value = "'191.79"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.91%  '

$ws.Range('D22').Value = "'10.12"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.85%  '

$ws.Range('D23').Value = "'6.030"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.71%  '

$ws.Range('E24').Value = '  +0.46%  '

$ws.Range('D25').Value = "'143.63"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.86%  '

$ws.Range('E26').Value = '  -1.80%  '

$ws.Range('D27').Value = "'7.226"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.37%  '

$ws.Range('D28').Value = "'16.03"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.25%  '

$ws.Range('D29').Value = "'1.473"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.80%  '

$ws.Range('D30').Value = "'0.05786"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.17%  '

$ws.Range('E31').Value = '  +0.12%  '

$ws.Range('D32').Value = "'3.565"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.81%  '

$ws.Range('D33').Value = "'3.281"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.28%  '

$ws.Range('D34').Value = "'1.601"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.83%  '

$ws.Range('D36').Value = "'0.9545"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.40%  '

$ws.Range('E37').Value = '  +0.37%  '

$ws.Range('D38').Value = "'0.5786"
$ws.Range('D38').Style = 'Normal'

$ws.Range('D39').Value = "'0.01604"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.01%  '

$ws.Range('E40').Value = '  -0.50%  '

$ws.Range('E41').Value = '  +0.81%  '

$ws.Range('E42').Value = '  +0.45%  '

$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = "'104.59"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.93%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = "'1.044.64"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.61%  '

$ws.Range('D45').Value = "'1.797.65"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.16%  '

$ws.Range('D46').Value = "'56.91"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.04%  '

$ws.Range('E47').Value = '  -1.40%  '

$ws.Range('D48').Value = "'1.008"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.86%  '

$ws.Range('D49').Value = "'0.4370"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.67%  '

$ws.Range('D50').Value = "'7.914"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.99%  '
